# Update TPM-derived values in the NATMI LR-pairs output for Ihh-Cdon.
# Only the numeric metrics derived from the new TPM values change; the
# categorical/label columns (A-D) and some counts stay the same except
# where noted below (K4/L4, K7/L7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (target cluster: ECs)
$ws.Range("M2").Value = 1.7455905
$ws.Range("N2").Value = 3.491181
$ws.Range("O2").Value = 0.06735438749894324
$ws.Range("P2").Value = 0.05664184489359518
$ws.Range("Q2").Value = 0.0271398592305
$ws.Range("R2").Value = 0.162839155383
$ws.Range("S2").Value = 0.06735438749894324
$ws.Range("T2").Value = 0.05664184489359518

# Row 3 (target cluster: FAPs)
$ws.Range("O3").Value = 0.3705234435972038
$ws.Range("P3").Value = 0.4673889601179179
$ws.Range("S3").Value = 0.3705234435972038
$ws.Range("T3").Value = 0.4673889601179179

# Row 4 (target cluster: Inflammatory-Mac)
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03863166666666667
$ws.Range("N4").Value = 0.115895
$ws.Range("O4").Value = 0.001490620077501955
$ws.Range("P4").Value = 0.001880311165173966
$ws.Range("Q4").Value = 0.0006006322761111111
$ws.Range("R4").Value = 0.005405690485
$ws.Range("S4").Value = 0.001490620077501955
$ws.Range("T4").Value = 0.001880311165173966

# Row 5 (target cluster: MuSCs)
$ws.Range("M5").Value = 14.3678575
$ws.Range("N5").Value = 28.735715
$ws.Range("O5").Value = 0.5543901857764452
$ws.Range("P5").Value = 0.4662158484296736
$ws.Range("Q5").Value = 0.2233866591241667
$ws.Range("R5").Value = 1.340319954745
$ws.Range("S5").Value = 0.5543901857764452
$ws.Range("T5").Value = 0.4662158484296736

# Row 6 (target cluster: Neutrophils)
$ws.Range("M6").Value = 0.06980833333333333
$ws.Range("N6").Value = 0.209425
$ws.Range("O6").Value = 0.002693585657110719
$ws.Range("P6").Value = 0.003397766648833495
$ws.Range("Q6").Value = 0.001085356697222222
$ws.Range("R6").Value = 0.009768210275
$ws.Range("S6").Value = 0.002693585657110719
$ws.Range("T6").Value = 0.003397766648833495

# Row 7 (target cluster: Resolving-Mac)
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.09194599999999999
$ws.Range("N7").Value = 0.275838
$ws.Range("O7").Value = 0.003547777392795065
$ws.Range("P7").Value = 0.004475268744805699
$ws.Range("Q7").Value = 0.001429545759333333
$ws.Range("R7").Value = 0.012865911834
$ws.Range("S7").Value = 0.003547777392795065
$ws.Range("T7").Value = 0.004475268744805699
